# Apply updated dSF (column F) values for specific rows, per commit:
# "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F4").Value = -2
$ws.Range("F7").Value = -5
$ws.Range("F14").Value = 4
$ws.Range("F21").Value = -7
$ws.Range("F22").Value = -8
$ws.Range("F24").Value = -2
$ws.Range("F25").Value = 1
$ws.Range("F26").Value = -2
$ws.Range("F27").Value = 2
$ws.Range("F28").Value = -1
$ws.Range("F31").Value = -2
